$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(96).Insert()

$ws.Range("A96").Value = 11
$ws.Range("B96").Value = "Vega Monumental Concepción"
$ws.Range("C96").Value = "Bíobío"
$ws.Range("D96").Value = 44778
$ws.Range("E96").Value = 8
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100108
$ws.Range("H96").Value = "Tropicales y subtropicales"
$ws.Range("I96").Value = 100108005
$ws.Range("J96").Value = "Piña"
$ws.Range("K96").Value = "Caramelo"
$ws.Range("L96").Value = "Tercera"
$ws.Range("M96").Value = 250
$ws.Range("N96").Value = 17000
$ws.Range("O96").Value = 18000
$ws.Range("P96").Value = 17600
$ws.Range("Q96").Value = "$/caja 16 unidades"
$ws.Range("R96").Value = "Ecuador"
$ws.Range("S96").Value = 1100
$ws.Range("T96").Value = 16
